$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 "REQUIRED"/"OPTIONAL" cells are re-shuffled by the merge:
#   A2 REQUIRED -> OPTIONAL
#   B2 REQUIRED -> REQUIRED (unchanged)
#   C2 REQUIRED -> OPTIONAL
#   D2 REQUIRED -> REQUIRED (unchanged)
#   E2 OPTIONAL -> OPTIONAL (unchanged)
#   F2 OPTIONAL -> OPTIONAL (unchanged)
$ws.Range("A2").Value = "OPTIONAL"
$ws.Range("B2").Value = "REQUIRED"
$ws.Range("C2").Value = "OPTIONAL"
$ws.Range("D2").Value = "REQUIRED"
$ws.Range("E2").Value = "OPTIONAL"
$ws.Range("F2").Value = "OPTIONAL"

# A2 previously carried a leftover one-off font (distinct from the rest of
# row 2); normalize it back to the same plain formatting as its neighbors
# so the redundant font/style entries drop out of the style table.
$ws.Range("A2").Font.Name = "Calibri"
$ws.Range("A2").Font.Size = 12
$ws.Range("A2").Font.Color = 0
$ws.Range("A2").Font.Bold = $false

# Column A width tweak (target raw width 19.6; Excel's ColumnWidth setter
# snaps to whole-pixel increments via the sheet's Normal-style font metrics,
# so 18.8 "characters" is the closest input that lands on that pixel grid).
$ws.Columns.Item(1).ColumnWidth = 18.8

# Move the active selection.
$ws.Range("G5").Select()
